$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41/42 swap: dogwifhat <-> Kaspa (coin name, link, price, change)
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"

# Price (column D) updates -- forced to remain text values
$ws.Range("D2").Value = "63.853.14"
$ws.Range("D3").Value = "3.195.18"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "3.191.95"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "3.717.16"
$ws.Range("D17").Value = "3.187.70"
$ws.Range("D18").Value = "63.839.03"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "0.0₃0744"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0396"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "397.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "2.809.92"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.255"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "129.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.72"
$ws.Range("D50").Style = "Normal"

# Volume(1h) percent change (column E) updates
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").Value = "  -2.67%  "
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("E43").Value = "  -5.96%  "
$ws.Range("E44").Value = "  -7.55%  "
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -0.58%  "
